# Scheduled market-data refresh: updates the computed price / profit
# columns (H:N) for a batch of rows across every job sheet (ALC, ARM,
# BSM, CRP, CUL, GSM, LTW, WVR) with freshly pulled Market Board values.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1521.3636
$ws.Range("I17").Value = 1388.6364
$ws.Range("J17").Value = 1587.7273
$ws.Range("K17").Value = 4165.9092
$ws.Range("L17").Value = 4763.1819
$ws.Range("M17").Value = -3997.9092
$ws.Range("N17").Value = -5099.1819
$ws.Range("H32").Value = 343.36365
$ws.Range("I32").Value = 347.75
$ws.Range("J32").Value = 331.66666
$ws.Range("K32").Value = 347.75
$ws.Range("L32").Value = 331.66666
$ws.Range("M32").Value = -21.75
$ws.Range("N32").Value = -983.66666
$ws.Range("H98").Value = 2884.4814
$ws.Range("I98").Value = 1786.7084
$ws.Range("K98").Value = 1786.7084
$ws.Range("M98").Value = -288.7084
$ws.Range("H112").Value = 10205655
$ws.Range("J112").Value = 1605.9579
$ws.Range("L112").Value = 4817.8737
$ws.Range("N112").Value = -7033.8737
$ws.Range("H122").Value = 2884.4814
$ws.Range("I122").Value = 1786.7084
$ws.Range("K122").Value = 5360.1252
$ws.Range("M122").Value = -2910.1252
$ws.Range("H129").Value = 36115.297
$ws.Range("J129").Value = 38973.36
$ws.Range("L129").Value = 116920.08
$ws.Range("N129").Value = -126920.08
$ws.Range("H132").Value = 100888.734
$ws.Range("I132").Value = 131208
$ws.Range("K132").Value = 393624
$ws.Range("M132").Value = -391094
$ws.Range("H133").Value = 30188.889
$ws.Range("J133").Value = 30188.889
$ws.Range("L133").Value = 30188.889
$ws.Range("N133").Value = -40308.889
$ws.Range("H138").Value = 2587.43
$ws.Range("I138").Value = 1079.4482
$ws.Range("J138").Value = 3203.3662
$ws.Range("K138").Value = 3238.3446
$ws.Range("L138").Value = 9610.098599999999
$ws.Range("M138").Value = 1901.6554
$ws.Range("N138").Value = -19890.0986

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2429.6956
$ws.Range("I61").Value = 1617.4546
$ws.Range("J61").Value = 3174.25
$ws.Range("K61").Value = 1617.4546
$ws.Range("L61").Value = 3174.25
$ws.Range("M61").Value = -1405.4546
$ws.Range("N61").Value = -3598.25
$ws.Range("H63").Value = 5330531
$ws.Range("I63").Value = 9895872
$ws.Range("K63").Value = 9895872
$ws.Range("M63").Value = -9895186
$ws.Range("H66").Value = 5330531
$ws.Range("I66").Value = 9895872
$ws.Range("K66").Value = 49479360
$ws.Range("M66").Value = -49475928
$ws.Range("H74").Value = 3475.4285
$ws.Range("I74").Value = 3567.7576
$ws.Range("K74").Value = 3567.7576
$ws.Range("M74").Value = -2693.7576
$ws.Range("H77").Value = 3475.4285
$ws.Range("I77").Value = 3567.7576
$ws.Range("K77").Value = 17838.788
$ws.Range("M77").Value = -13470.788
$ws.Range("H132").Value = 2950.0244
$ws.Range("I132").Value = 2163.3215
$ws.Range("J132").Value = 4644.4614
$ws.Range("K132").Value = 6489.9645
$ws.Range("L132").Value = 13933.3842
$ws.Range("M132").Value = -3959.9645
$ws.Range("N132").Value = -18993.3842
$ws.Range("H136").Value = 2429.6956
$ws.Range("I136").Value = 1617.4546
$ws.Range("J136").Value = 3174.25
$ws.Range("K136").Value = 4852.3638
$ws.Range("L136").Value = 9522.75
$ws.Range("M136").Value = -2302.3638
$ws.Range("N136").Value = -14622.75
$ws.Range("H137").Value = 51363.6
$ws.Range("J137").Value = 51363.6
$ws.Range("L137").Value = 51363.6
$ws.Range("N137").Value = -61563.6

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H137").Value = 32795
$ws.Range("J137").Value = 32795
$ws.Range("L137").Value = 32795
$ws.Range("N137").Value = -42995

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2080.59
$ws.Range("I58").Value = 1753.0182
$ws.Range("J58").Value = 5083.3335
$ws.Range("K58").Value = 1753.0182
$ws.Range("L58").Value = 5083.3335
$ws.Range("M58").Value = -1550.0182
$ws.Range("N58").Value = -5489.3335
$ws.Range("H132").Value = 6116.778
$ws.Range("I132").Value = 6803.4287
$ws.Range("K132").Value = 20410.2861
$ws.Range("M132").Value = -17880.2861
$ws.Range("H136").Value = 2080.59
$ws.Range("I136").Value = 1753.0182
$ws.Range("J136").Value = 5083.3335
$ws.Range("K136").Value = 5259.054599999999
$ws.Range("L136").Value = 15250.0005
$ws.Range("M136").Value = -2709.054599999999
$ws.Range("N136").Value = -20350.0005

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1912.0555
$ws.Range("I5").Value = 432.07693
$ws.Range("K5").Value = 1296.23079
$ws.Range("M5").Value = -1184.23079
$ws.Range("H46").Value = 1498.5834
$ws.Range("I46").Value = 503
$ws.Range("J46").Value = 1589.091
$ws.Range("K46").Value = 1509
$ws.Range("L46").Value = 4767.272999999999
$ws.Range("M46").Value = -1418
$ws.Range("N46").Value = -4949.272999999999
$ws.Range("H86").Value = 6746.923
$ws.Range("I86").Value = 7274.75
$ws.Range("K86").Value = 21824.25
$ws.Range("M86").Value = -20638.25
$ws.Range("H89").Value = 6746.923
$ws.Range("I89").Value = 7274.75
$ws.Range("K89").Value = 65472.75
$ws.Range("M89").Value = -59544.75
$ws.Range("H113").Value = 645.8125
$ws.Range("I113").Value = 625.2593000000001
$ws.Range("J113").Value = 756.8
$ws.Range("K113").Value = 1875.7779
$ws.Range("L113").Value = 2270.4
$ws.Range("M113").Value = 294.2221
$ws.Range("N113").Value = -6610.4
$ws.Range("H125").Value = 10000
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 10000
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 30000
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -39840
$ws.Range("H131").Value = 10639259
$ws.Range("I131").Value = 100002240
$ws.Range("J131").Value = 809.0238000000001
$ws.Range("K131").Value = 300006720
$ws.Range("L131").Value = 2427.0714
$ws.Range("M131").Value = -300001680
$ws.Range("N131").Value = -12507.0714
$ws.Range("H132").Value = 2401.3618
$ws.Range("I132").Value = 918
$ws.Range("J132").Value = 3166.9678
$ws.Range("K132").Value = 8262
$ws.Range("L132").Value = 28502.7102
$ws.Range("M132").Value = -5732
$ws.Range("N132").Value = -33562.7102
$ws.Range("H135").Value = 1912.0555
$ws.Range("I135").Value = 432.07693
$ws.Range("K135").Value = 3888.69237
$ws.Range("M135").Value = -1353.69237

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 620.3333
$ws.Range("I107").Value = 396.33334
$ws.Range("J107").Value = 844.3333
$ws.Range("K107").Value = 396.33334
$ws.Range("L107").Value = 844.3333
$ws.Range("M107").Value = 1523.66666
$ws.Range("N107").Value = -4684.3333
$ws.Range("H113").Value = 1230.2667
$ws.Range("I113").Value = 1211
$ws.Range("K113").Value = 1211
$ws.Range("M113").Value = 959
$ws.Range("H132").Value = 2806.1785
$ws.Range("I132").Value = 1319.75
$ws.Range("J132").Value = 4788.0835
$ws.Range("K132").Value = 3959.25
$ws.Range("L132").Value = 14364.2505
$ws.Range("M132").Value = -1429.25
$ws.Range("N132").Value = -19424.2505
$ws.Range("H134").Value = 49281.523
$ws.Range("J134").Value = 50780.8
$ws.Range("L134").Value = 152342.4
$ws.Range("N134").Value = -157412.4
$ws.Range("H137").Value = 72533.11
$ws.Range("J137").Value = 72533.11
$ws.Range("L137").Value = 72533.11
$ws.Range("N137").Value = -82733.11

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 3332.2222
$ws.Range("I93").Value = 1996.6666
$ws.Range("K93").Value = 1996.6666
$ws.Range("M93").Value = -748.6666
$ws.Range("H132").Value = 4007.2415
$ws.Range("I132").Value = 1655.7391
$ws.Range("J132").Value = 5552.514
$ws.Range("K132").Value = 4967.2173
$ws.Range("L132").Value = 16657.542
$ws.Range("M132").Value = -2437.2173
$ws.Range("N132").Value = -21717.542
$ws.Range("H136").Value = 3659.1143
$ws.Range("I136").Value = 1636.2084
$ws.Range("J136").Value = 8072.727
$ws.Range("K136").Value = 4908.6252
$ws.Range("L136").Value = 24218.181
$ws.Range("M136").Value = -2358.6252
$ws.Range("N136").Value = -29318.181

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2906.3684
$ws.Range("I136").Value = 1380.2307
$ws.Range("K136").Value = 4140.6921
$ws.Range("M136").Value = -1590.6921

